$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '19.940.28'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -5.46%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.414.13'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -6.24%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.56%  '
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '275.92'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3666'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -5.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3097'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '39.76'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -6.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.034'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06524'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -6.69%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.483'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.56'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.187'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.415.37'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -6.48%  '
$ws.Range('E17').Value = '  -4.91%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.05659'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -13.89%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '71.11'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -13.31%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.609'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.44%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.70'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.94%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.89'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.236'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.78%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '19.959.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.257'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.74%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '132.90'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -10.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.24'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.33%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.571.29'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -6.80%  '
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.893'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -18.64%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.272'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -10.95%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8185'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -14.47%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07683'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.489'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '8.311'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.914'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05781'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9972'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02054'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.41%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '10.44'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -7.32%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1884'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.52%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.092'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.53%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '12.41'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.82%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5305'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.66%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.536'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.67%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5177'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '115.06'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.765'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.54%  '
$ws.Range('E50').Value = '  -9.09%  '
$ws.Range('E51').Value = '  -0.63%  '
